$d = $word.ActiveDocument

# Change 1: "University employees who are students" -> "University employees who have been students"
$d.Content.Find.Execute("University employees who are students", $true, $false, $false, $false, $false, $true, 1, $false, "University employees who have been students", 2)

# Change 2: "apply through the Global Engagement Office." -> "apply through Global Engagement."
$d.Content.Find.Execute("apply through the Global Engagement Office.", $true, $false, $false, $false, $false, $true, 1, $false, "apply through Global Engagement.", 2)

# Change 3: explicitly set page orientation to portrait so w:orient="portrait" is emitted on pgSz
$d.PageSetup.Orientation = 0
